$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.1724137931034483
$ws.Range("C2").Value = 0.6206896551724138
$ws.Range("P2").Value = 0.1379310344827586
$ws.Range("S2").Value = 0.06896551724137931

# Row 3
$ws.Range("P3").Value = 0.8235294117647058
$ws.Range("S3").Value = 0.1764705882352941

# Row 4
$ws.Range("P4").Value = 1

# Row 6
$ws.Range("B6").Value = 0.09523809523809523
$ws.Range("F6").Value = 0.09523809523809523
$ws.Range("J6").Value = 0.04761904761904762
$ws.Range("Q6").Value = 0.09523809523809523
$ws.Range("R6").Value = 0.04761904761904762
$ws.Range("S6").Value = 0.6190476190476191

# Row 7
$ws.Range("B7").Value = 0.04166666666666666
$ws.Range("F7").Value = 0.04166666666666666
$ws.Range("J7").Value = 0.08333333333333333
$ws.Range("Q7").Value = 0.08333333333333333
$ws.Range("S7").Value = 0.75

# Row 8
$ws.Range("B8").Value = 0.1052631578947368
$ws.Range("F8").Value = 0.04210526315789474
$ws.Range("J8").Value = 0.03157894736842105
$ws.Range("O8").Value = 0.02105263157894737
$ws.Range("Q8").Value = 0.1473684210526316
$ws.Range("R8").Value = 0.09473684210526316
$ws.Range("S8").Value = 0.5578947368421052

# Row 9
$ws.Range("B9").Value = 0.04761904761904762
$ws.Range("D9").Value = 0.04761904761904762
$ws.Range("F9").Value = 0.04761904761904762
$ws.Range("J9").Value = 0.04761904761904762
$ws.Range("Q9").Value = 0.04761904761904762
$ws.Range("S9").Value = 0.7619047619047619

# Row 10
$ws.Range("B10").Value = 0.08849557522123894
$ws.Range("F10").Value = 0.08849557522123894
$ws.Range("J10").Value = 0.01769911504424779
$ws.Range("O10").Value = 0.008849557522123894
$ws.Range("Q10").Value = 0.1238938053097345
$ws.Range("R10").Value = 0.07964601769911504
$ws.Range("S10").Value = 0.5929203539823009

# Row 11
$ws.Range("G11").Value = 0.225
$ws.Range("J11").Value = 0.05
$ws.Range("K11").Value = 0.275
$ws.Range("L11").Value = 0.425
$ws.Range("S11").Value = 0.025

# Row 12
$ws.Range("G12").Value = 0.6470588235294118
$ws.Range("J12").Value = 0.2941176470588235
$ws.Range("S12").Value = 0.05882352941176471

# Row 13
$ws.Range("G13").Value = 0.5714285714285714
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.1428571428571428

# Row 15
$ws.Range("H15").Value = 0.4736842105263158
$ws.Range("J15").Value = 0.2105263157894737
$ws.Range("M15").Value = 0.05263157894736842
$ws.Range("O15").Value = 0.05263157894736842
$ws.Range("S15").Value = 0.2105263157894737

# Row 16
$ws.Range("F16").Value = 0.05555555555555555
$ws.Range("I16").Value = 0.1111111111111111
$ws.Range("J16").Value = 0.2777777777777778
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("M16").Value = 0.05555555555555555
$ws.Range("S16").Value = 0.3888888888888889

# Row 17
$ws.Range("H17").Value = 0.2424242424242424
$ws.Range("I17").Value = 0.1212121212121212
$ws.Range("J17").Value = 0.2727272727272727
$ws.Range("K17").Value = 0.06060606060606061
$ws.Range("M17").Value = 0.0303030303030303
$ws.Range("O17").Value = 0.06060606060606061
$ws.Range("S17").Value = 0.2121212121212121

# Row 18
$ws.Range("H18").Value = 0.2105263157894737
$ws.Range("I18").Value = 0.1052631578947368
$ws.Range("J18").Value = 0.4736842105263158
$ws.Range("K18").Value = 0.1052631578947368
$ws.Range("S18").Value = 0.1052631578947368

# Row 19
$ws.Range("F19").Value = 0.004310344827586207
$ws.Range("H19").Value = 0.3275862068965517
$ws.Range("I19").Value = 0.05603448275862069
$ws.Range("J19").Value = 0.293103448275862
$ws.Range("K19").Value = 0.09913793103448276
$ws.Range("M19").Value = 0.01724137931034483
$ws.Range("O19").Value = 0.03879310344827586
$ws.Range("S19").Value = 0.1637931034482759
